$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New Q04 quiz column (K) -------------------------------------------
# Header
$ws.Range("K1").Value = "Q04"

# Row 2
$ws.Range("H2").Value = 5.5
$ws.Range("K2").Formula = "=(11/17)*10"

# Row 3
$ws.Range("H3").Value = 9
$ws.Range("K3").Formula = "=(9/17)*10"

# Row 4
$ws.Range("H4").Value = 9.5
$ws.Range("K4").Formula = "=(15/17)*10"

# Row 5
$ws.Range("H5").Value = 6
$ws.Range("K5").Formula = "=(10/17)*10"

# Row 6
$ws.Range("H6").Value = 8.5
$ws.Range("K6").Formula = "=(10/17)*10"

# Row 7
$ws.Range("H7").Value = 9
$ws.Range("K7").Formula = "=(12/17)*10"

# Row 8 (plain value, not a formula)
$ws.Range("K8").Value = 0

# Row 9
$ws.Range("H9").Value = 9.5
$ws.Range("K9").Formula = "=(11/17)*10"

# Row 10
$ws.Range("H10").Value = 8
$ws.Range("K10").Formula = "=(8/17)*10"

# Row 11
$ws.Range("H11").Value = 6
$ws.Range("K11").Formula = "=7/17*10"
$ws.Range("K11").HorizontalAlignment = -4108

# Row 12
$ws.Range("H12").Value = 6.5
$ws.Range("K12").Formula = "=10/17*10"
$ws.Range("K12").HorizontalAlignment = -4108

# Row 13
$ws.Range("H13").Value = 9
$ws.Range("K13").Formula = "=11/17*10"
$ws.Range("K13").HorizontalAlignment = -4108

# Row 14
$ws.Range("H14").Value = 5.5
$ws.Range("K14").Formula = "=7/17*10"
$ws.Range("K14").HorizontalAlignment = -4108

# Row 15
$ws.Range("H15").Value = 6.5
$ws.Range("K15").Formula = "=10/17*10"
$ws.Range("K15").HorizontalAlignment = -4108

# Row 16
$ws.Range("H16").Value = 7.5
$ws.Range("K16").Formula = "=5/17*10"
$ws.Range("K16").HorizontalAlignment = -4108

# Row 17 totals - K17 only needs the matching centered format (no value)
$ws.Range("K17").HorizontalAlignment = -4108

# --- Selection cursor, as recorded in the saved view --------------------
[void]$ws.Range("H22").Select()
